$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A59").Value = "GRT-USD"
